# Daily attendance processing - 2026-01-01 08:40:10
# Rotates the "Recorded By" list in column G so the last recorder
# (the most recent one to touch the record) moves to the front of the list.
# Cells with only a single recorder are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ", "
    $n = $parts.Length

    if ($n -gt 1) {
        $last = $parts[$n - 1]
        $rest = $parts[0..($n - 2)]
        $newParts = @($last) + $rest
        $newVal = $newParts -join ", "
        $cell.Value2 = $newVal
    }
}
